$d = $word.ActiveDocument

# 1. Remove the first paragraph ("-create grocery table") entirely -- it
#    was replaced by the (now-renumbered) paragraphs below.
$d.Paragraphs(1).Range.Delete()

# 2. Append " (Name)" as its own trailing run to each remaining paragraph.
$names = @("Carson", "Brooke", "Brooke", "Brooke", "Carson")

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $name = $names[$i - 1]
    $suffix = " (" + $name + ")"

    # Insertion point: right before the paragraph mark.
    $insertAt = $para.Range.End - 1
    $ins = $d.Range($insertAt, $insertAt)
    $ins.InsertAfter($suffix)

    # Force the newly inserted text into its own run (a distinct <w:r>,
    # matching the target markup) by toggling a character-formatting
    # property just across that span, then reverting it.
    $newRange = $d.Range($insertAt, $insertAt + $suffix.Length)
    $newRange.Bold = 1
    $newRange.Bold = 0
}

# 3. Relocate the hidden "_GoBack" bookmark so it ends up collapsed at the
#    very end of the final paragraph (after " (Carson)"), mirroring how
#    Word itself drags that bookmark along with the last edit.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # No pre-existing _GoBack bookmark -- nothing to relocate away from.
}

$last = $d.Paragraphs($d.Paragraphs.Count)
$endPos = $last.Range.End - 1

# A bookmark collapsed exactly on the paragraph-mark boundary can't be
# placed directly, so temporarily park a throwaway character there to
# give the insertion point real neighbours, bookmark it, then clean up.
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Range.Text + "]")
}
